$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 23:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1641220
$ws.Range("C4").Value = 20318
$ws.Range("D4").Value = 395551
$ws.Range("E4").Value = 1148177
$ws.Range("G4").Value = 1138
$ws.Range("H4").Value = 97492

# Row 17 - Canada
$ws.Range("B17").Value = 82420
$ws.Range("C17").Value = 1096
$ws.Range("D17").Value = 42481
$ws.Range("E17").Value = 33694

# Row 63 - Oman
$ws.Range("E63").Value = 4939
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = 34

# Row 84 - Costa de Marfil
$ws.Range("B84").Value = 2341
$ws.Range("C84").Value = 40
$ws.Range("D84").Value = 1146
$ws.Range("E84").Value = 1166

# Row 188 - Botsuana
$ws.Range("B188").Value = 30
$ws.Range("C188").Value = 1
$ws.Range("E188").Value = 10

$wb.Save()
